# Add a new "sortQuoteList" test-data block to the DashboardPageData sheet,
# mirroring the existing blocks (title row + header row + Y/N data rows)
# that already live at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DashboardPageData")

# --- Title row (A55:B55) -------------------------------------------------
# Copy values + formatting from the block above (row 50), then overwrite
# the label with the new test name "sortQuoteList".
$ws.Range("A50:B50").Copy()
$ws.Range("A55").PasteSpecial(-4163)   # xlPasteValuesAndNumberFormats
$ws.Range("A50:B50").Copy()
$ws.Range("A55").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A55").Value = "sortQuoteList"

# --- Header row (A56:D56) -------------------------------------------------
$ws.Range("A51:D51").Copy()
$ws.Range("A56").PasteSpecial(-4163)
$ws.Range("A51:D51").Copy()
$ws.Range("A56").PasteSpecial(-4122)

# --- Data row, runMode = Y (A57:D57) --------------------------------------
$ws.Range("A52:D52").Copy()
$ws.Range("A57").PasteSpecial(-4163)
$ws.Range("A52:D52").Copy()
$ws.Range("A57").PasteSpecial(-4122)

# --- Data row, runMode = N (A58:D58) --------------------------------------
# No "N" row existed in this particular block before, so pull the matching
# pattern (same brokerId/agentId/agencyOfficeId trio used elsewhere in the
# sheet for the N case) from row 47.
$ws.Range("A47:D47").Copy()
$ws.Range("A58").PasteSpecial(-4163)
$ws.Range("A47:D47").Copy()
$ws.Range("A58").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Restore the view so it's scrolled/selected near the new block -------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 32
$ws.Range("C67").Select() | Out-Null
